$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the latest round-of-16 match scores (previously blank cells)
$ws.Range("F47").Value = 0   # Serbia
$ws.Range("G47").Value = 2   # Brazil

$ws.Range("F48").Value = 2   # Switzerland
$ws.Range("G48").Value = 2   # Costa Rica

$ws.Range("F49").Value = 2   # Korea Republic
$ws.Range("G49").Value = 0   # Germany

$ws.Range("F50").Value = 0   # Mexico
$ws.Range("G50").Value = 3   # Sweden

# Leave the view scrolled/selected where the user was last working
$ws.Range("F53").Select()
